# Weekly update: insert 3 new rows (Especial/Primera/Segunda) of fresh daily
# data at the top of the price history table for
# "Terminal La Palmera de La Serena - Mango", pushing all existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 357 (existing rows 357:413 shift to 360:416).
$ws.Range("A357:A359").EntireRow.Insert()

# Populate the 3 new rows with this week's record (the columns that stay
# constant for every row of this sheet are copied from the template/pattern
# used throughout the table; only the day's price/volume/origin data change).
$qualities = @("Especial", "Primera", "Segunda")
for ($i = 0; $i -lt 3; $i++) {
    $r = 357 + $i
    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44476
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108002
    $ws.Cells.Item($r, 10).Value = "Mango"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $qualities[$i]
    $ws.Cells.Item($r, 13).Value = 512
    $ws.Cells.Item($r, 14).Value = 7000
    $ws.Cells.Item($r, 15).Value = 7500
    $ws.Cells.Item($r, 16).Value = 7250
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 4 kilos"
    $ws.Cells.Item($r, 18).Value = "Perú"
    $ws.Cells.Item($r, 19).Value = 1812
    $ws.Cells.Item($r, 20).Value = 4
}
